# Updated cryptos list (refresh of price/volume data scraped into the sheet).
# Mirrors the GitHub Actions commit "Updated cryptos list ... with GitHub Actions".
#
# Note: several "Price" (column D) values look like plain numbers (e.g. "8.40",
# "77.00", "11.90") but are stored as TEXT in the workbook (European-style
# thousand separators like "42.713.69" make the column text-typed overall).
# Assigning such a string straight to .Value would make Excel coerce it to a
# number and silently drop the significant trailing zero (e.g. "8.40" -> 8.4).
# Prefixing with a leading apostrophe forces Excel to keep it as text, exactly
# like a user typing '8.40 into the cell would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.713.69"
$ws.Range("E2").Value = "  -1.42%  "

$ws.Range("D3").Value = "2.348.74"
$ws.Range("E3").Value = "  -1.81%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'321.14"
$ws.Range("E5").Value = "  -0.76%  "

$ws.Range("D6").Value = "'105.15"
$ws.Range("E6").Value = "  -0.54%  "

$ws.Range("D7").Value = "'0.636"
$ws.Range("E7").Value = "  -2.88%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").Value = "'0.615"
$ws.Range("E9").Value = "  -6.09%  "

$ws.Range("D10").Value = "'41.17"
$ws.Range("E10").Value = "  -1.92%  "

$ws.Range("D11").Value = "'0.0922"
$ws.Range("E11").Value = "  -2.51%  "

$ws.Range("D12").Value = "'8.40"
$ws.Range("E12").Value = "  -2.17%  "

$ws.Range("D13").Value = "'0.998"
$ws.Range("E13").Value = "  -1.90%  "

$ws.Range("E14").Value = "  -0.49%  "

$ws.Range("D15").Value = "'15.91"
$ws.Range("E15").Value = "  -8.09%  "

$ws.Range("D16").Value = "2.704.10"
$ws.Range("E16").Value = "  -1.88%  "

$ws.Range("D17").Value = "2.325.32"
$ws.Range("E17").Value = "  -2.82%  "

$ws.Range("D18").Value = "42.649.00"
$ws.Range("E18").Value = "  -1.53%  "

$ws.Range("D19").Value = "'7.72"
$ws.Range("E19").Value = "  +4.57%  "

$ws.Range("E20").Value = "  -2.90%  "

$ws.Range("D21").Value = "'77.00"
$ws.Range("E21").Value = "  +1.62%  "

$ws.Range("D22").Value = "'3.63"
$ws.Range("E22").Value = "  +5.13%  "

$ws.Range("D23").Value = "'258.41"
$ws.Range("E23").Value = "  -3.63%  "

$ws.Range("D24").Value = "'2.30"
$ws.Range("E24").Value = "  -5.08%  "

$ws.Range("D25").Value = "'9.41"
$ws.Range("E25").Value = "  -3.71%  "

$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("D27").Value = "'11.36"
$ws.Range("E27").Value = "  -4.51%  "

$ws.Range("D28").Value = "'22.96"
$ws.Range("E28").Value = "  +0.46%  "

$ws.Range("E29").Value = "  +0.38%  "

$ws.Range("D30").Value = "'174.77"
$ws.Range("E30").Value = "  -1.54%  "

$ws.Range("D31").Value = "'36.25"
$ws.Range("E31").Value = "  -4.35%  "

$ws.Range("D32").Value = "'0.0887"
$ws.Range("E32").Value = "  -4.55%  "

$ws.Range("D33").Value = "'6.09"
$ws.Range("E33").Value = "  +3.24%  "

$ws.Range("D34").Value = "'2.95"
$ws.Range("E34").Value = "  -8.21%  "

$ws.Range("D35").Value = "'0.122"
$ws.Range("E35").Value = "  +11.35%  "

$ws.Range("E36").Value = "  -3.33%  "

$ws.Range("D37").Value = "'4.59"
$ws.Range("E37").Value = "  -5.63%  "

$ws.Range("D38").Value = "'0.0359"
$ws.Range("E38").Value = "  -2.53%  "

$ws.Range("D39").Value = "'3.76"
$ws.Range("E39").Value = "  -9.30%  "

$ws.Range("D40").Value = "'2.69"
$ws.Range("E40").Value = "  -4.83%  "

$ws.Range("D41").Value = "'71.21"
$ws.Range("E41").Value = "  +2.50%  "

$ws.Range("D42").Value = "'0.234"
$ws.Range("E42").Value = "  +0.36%  "

$ws.Range("D43").Value = "'1.46"
$ws.Range("E43").Value = "  -7.57%  "

$ws.Range("E44").Value = "  -0.18%  "

$ws.Range("D45").Value = "'114.82"
$ws.Range("E45").Value = "  -8.94%  "

# Rows 46/47 swap places (Celestia now ranks above BitcoinSV) with refreshed
# price/volume figures.
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").Value = "'11.90"
$ws.Range("E46").Value = "  -5.83%  "

$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").Value = "'86.52"
$ws.Range("E47").Value = "  +1.38%  "

$ws.Range("D48").Value = "'5.49"
$ws.Range("E48").Value = "  -3.07%  "

$ws.Range("D49").Value = "'9.12"
$ws.Range("E49").Value = "  -5.36%  "

$ws.Range("D50").Value = "'73.42"
$ws.Range("E50").Value = "  +0.84%  "

$ws.Range("E51").Value = "  -1.40%  "
